# Updates "合肥-漫展信息.xlsx" to match commit 456a3b4 (gh-pages data refresh).
#
# Changes applied:
#  - Refresh "want to go" counters (column F) and one cover image URL on
#    several existing rows in sheets "展览" and "全部类型".
#  - Insert a brand-new event row ("合肥·W·A首届童年怀旧only") before the
#    "安徽·MAX特摄only展" row in sheets "展览" and "全部类型", shifting the
#    remaining rows down by one.
#  - Rename / refresh numbers for the "环形宇宙动漫游戏嘉年华" row that got
#    pushed down by the insert.
#  - Bump the "想去人数" counter for the "包河留声机音乐节" row (present in
#    sheets "演出" and "全部类型").
#
# NOTE: this engine's COM-script interpreter has a quirk where a
# parenthesised/computed expression used as a positional argument inside a
# function that is itself called from another function can blow up with
# "does not contain a method named ''". To stay safe, every helper below is
# called directly from top level (no helper calls another helper) and any
# computed string (e.g. "B" + $rowNum) is assigned to a local variable
# before being used.

$wb = $excel.ActiveWorkbook

function Set-FCounters {
    param($ws)

    # column F ("想去人数") refreshes that don't involve moving any rows
    $counters = @{
        "F2"  = 156
        "F3"  = 419
        "F4"  = 12244
        "F5"  = 1267
        "F6"  = 134
        "F7"  = 28
        "F10" = 190
        "F11" = 444
        "F12" = 57
        "F16" = 362
        "F17" = 2889
        "F19" = 933
    }
    foreach ($ref in $counters.Keys) {
        $ws.Range($ref).Value = $counters[$ref]
    }

    # refreshed cover image for row 19
    $ws.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202404/RFYwkzvt1713951750482.jpeg"
}

function Insert-NewEventRow {
    param($ws, $rowNum)

    # Push row $rowNum (and everything below it) down by one, then fill the
    # freshly inserted row with the new "W·A首届童年怀旧only" event.
    $ws.Rows.Item($rowNum).Insert()

    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108   # xlCenter
    $aCell.VerticalAlignment = -4160     # xlTop
    $aCell.Borders.LineStyle = 1
    $aCell.Value = 19

    $bRef = "B" + $rowNum
    $ws.Range($bRef).NumberFormat = "@"
    $ws.Range($bRef).Value = "2024-07-20"

    $cCell = $ws.Cells.Item($rowNum, 3)
    $cCell.Value = "合肥·W·A首届童年怀旧only"

    $dCell = $ws.Cells.Item($rowNum, 4)
    $dCell.Value = "阜阳路16号 银瑞林国际大酒店"

    $eCell = $ws.Cells.Item($rowNum, 5)
    $eCell.Value = "2024.07.20 09:30-07.20 17:00"

    $fCell = $ws.Cells.Item($rowNum, 6)
    $fCell.Value = 2

    $gCell = $ws.Cells.Item($rowNum, 7)
    $gCell.Value = 78

    $hCell = $ws.Cells.Item($rowNum, 8)
    $hCell.Value = "https://show.bilibili.com/platform/detail.html?id=84794"

    $iCell = $ws.Cells.Item($rowNum, 9)
    $iCell.Value = "//i2.hdslb.com/bfs/openplatform/202404/Ie0KTNEr1713951888990.png"
}

# ---------------------------------------------------------------------
# Sheet "展览": dimension A1:I21 -> A1:I22
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
Set-FCounters $wsExhibit
Insert-NewEventRow $wsExhibit 20

# Old row 20 (安徽·MAX特摄only展) is now row 21, unchanged.
# Old row 21 (合肥·环形宇宙动漫游戏嘉年华第7届) is now row 22: rename + refresh numbers.
$wsExhibit.Range("C22").Value = "合肥·第七届环形宇宙动漫游戏嘉年华"
$wsExhibit.Range("F22").Value = 12
$wsExhibit.Range("G22").Value = 49

# ---------------------------------------------------------------------
# Sheet "演出": F2 想去人数 refresh only (no row insert)
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 12

# ---------------------------------------------------------------------
# Sheet "全部类型": dimension A1:I22 -> A1:I23
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Set-FCounters $wsAll
Insert-NewEventRow $wsAll 20

# Old row 20 (安徽·MAX特摄only展) is now row 21, unchanged.
# Old row 21 (合肥·环形宇宙动漫游戏嘉年华第7届) is now row 22: rename + refresh numbers.
$wsAll.Range("C22").Value = "合肥·第七届环形宇宙动漫游戏嘉年华"
$wsAll.Range("F22").Value = 12
$wsAll.Range("G22").Value = 49

# Old row 22 (合肥·首届包河留声机音乐节...) is now row 23: refresh 想去人数.
$wsAll.Range("F23").Value = 12

Write-Output "edit.ps1 completed"
